# The source deck was re-saved with Aspose.Slides for .NET 15.7.0.0, which
# stamps the presentation with its own custom "tags" metadata part and
# (for this particular slide) records an explicit "no fill" background
# override instead of inheriting the layout/master background.
#
# Reproduce both observable, content-level effects through the PowerPoint
# object model:

$p = $ppt.ActivePresentation

# 1) Presentation-level custom data tags (-> ppt/tags/tag1.xml +
#    <p:custDataLst><p:tags .../></p:custDataLst> in presentation.xml)
$p.Tags.Add("AS_NET", "4.0.30319.34209")
$p.Tags.Add("AS_OS", "Microsoft Windows NT 6.2.9200.0")
$p.Tags.Add("AS_RELEASE_DATE", "2015.08.28")
$p.Tags.Add("AS_TITLE", "Aspose.Slides for .NET 4.0")
$p.Tags.Add("AS_VERSION", "15.7.0.0")

# 2) Slide 1 gets its own explicit "no fill" background instead of
#    inheriting the master/layout background (-> <p:bg><p:bgPr><a:noFill/>
#    </p:bgPr></p:bg> in ppt/slides/slide1.xml)
$s = $p.Slides.Item(1)
$s.FollowMasterBackground = $false
$s.Background.Fill.Background()
$s.Background.Fill.Visible = $false
